$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns being touched stay as text (they hold
# strings like "297.77" or "-0.31%", not real numbers/percentages).
$ws.Range("D2:E50").NumberFormat = "@"

# Row 2 (BNB)
$ws.Range("D2").Value = "297.77"
$ws.Range("E2").Value = "-0.31%"

# Row 3 (OKB)
$ws.Range("E3").Value = "-0.35%"

# Row 4 (HuobiToken)
$ws.Range("E4").Value = "-0.99%"

# Row 5 (Cronos)
$ws.Range("D5").Value = "0.08010"
$ws.Range("E5").Value = "9.43%"

# Row 6 (FTXToken)
$ws.Range("D6").Value = "2.425"
$ws.Range("E6").Value = "33.89%"

# Row 7 (KuCoinToken)
$ws.Range("D7").Value = "7.796"

# Row 8 (was GateToken -> MXToken)
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9245"
$ws.Range("E8").Value = "-0.13%"

# Row 9 (was MXToken -> WazirX)
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1725"
$ws.Range("E9").Value = "2.83%"

# Row 10 (was WazirX -> LiechtensteinCryptoassetsExchange)
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.07322"
$ws.Range("E10").Value = "2.92%"

# Row 11 (was LiechtensteinCryptoassetsExchange -> MandalaExchangeToken)
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08825"
$ws.Range("E11").Value = "10.18%"

# Row 12 (was MandalaExchangeToken -> BitrueCoin)
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03034"
$ws.Range("E12").Value = "1.14%"

# Row 13 (was BitrueCoin -> BitMartToken)
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09997"
$ws.Range("E13").Value = "0.75%"

# Row 14 (was BitMartToken -> BitForexToken)
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001502"
$ws.Range("E14").Value = "0.29%"

# Row 15 (was BitForexToken -> TigerCash)
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005916"
$ws.Range("E15").Value = "-4.43%"

# Row 16 (was TigerCash -> LEO)
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.502"
$ws.Range("E16").Value = "1.32%"

# Row 17 (was LEO -> GateToken)
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.798"
$ws.Range("E17").Value = "1.53%"

# Row 18 (BTSEToken)
$ws.Range("D18").Value = "2.246"
$ws.Range("E18").Value = "1.06%"

# Row 19 (BitpandaEcosystemToken)
$ws.Range("E19").Value = "1.87%"

# Row 20 (ProBitToken)
$ws.Range("E20").Value = "1.47%"

# Row 21 (MCDex)
$ws.Range("E21").Value = "0.73%"

# Row 22 (ZBToken)
$ws.Range("D22").Value = "0.1616"
$ws.Range("E22").Value = "2.15%"

# Row 23 (CoinExToken)
$ws.Range("D23").Value = "0.04608"
$ws.Range("E23").Value = "-0.71%"

# Row 24 (BitKan)
$ws.Range("E24").Value = "2.44%"

# Row 25 (HotbitToken)
$ws.Range("D25").Value = "0.004431"
$ws.Range("E25").Value = "-6.31%"

# Row 26 (NitroEx)
$ws.Range("E26").Value = "-7.64%"

# Row 27 (UpBots)
$ws.Range("E27").Value = "82.78%"

# Row 39 (One)
$ws.Range("D39").Value = "0.01787"
$ws.Range("E39").Value = "4.43%"

# Row 40 (IDEX)
$ws.Range("D40").Value = "0.04464"
$ws.Range("E40").Value = "-0.06%"

# Row 41 (KickToken)
$ws.Range("D41").Value = "0.006974"
$ws.Range("E41").Value = "-1.46%"

# Row 42 (BKEXToken)
$ws.Range("D42").Value = "0.1345"
$ws.Range("E42").Value = "1.04%"

# Row 43 (CEJI)
$ws.Range("D43").Value = "0.002209"
$ws.Range("E43").Value = "3.82%"

# Row 44 (LocalTraders)
$ws.Range("D44").Value = "0.009801"
$ws.Range("E44").Value = "-6.06%"

# Row 45 (CoinLion)
$ws.Range("D45").Value = "0.00006570"
$ws.Range("E45").Value = "5.60%"

# Row 46 (Kangarootoken)
$ws.Range("E46").Value = "-0.03%"

# Row 48 (BOLO)
$ws.Range("E48").Value = "-57.27%"

# Row 49 (CryptobidCoin)
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "-0.03%"

# Row 50 (SpecialPowerGold)
$ws.Range("E50").Value = "0.04%"
